$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("展览")
$ws.Range("F3").Value = 541
$ws.Range("F6").Value = 2397
$ws.Range("F7").Value = 71
$ws.Range("F11").Value = 1506
$ws.Range("F13").Value = 597
$ws.Range("F14").Value = 739
$ws.Range("F15").Value = 1091
$ws.Range("F16").Value = 485
$ws.Range("F17").Value = 3404
$ws.Range("F19").Value = 634
$ws.Range("F20").Value = 3250
$ws.Range("F21").Value = 721
$ws.Range("F22").Value = 601
$ws.Range("F23").Value = 15
$ws.Range("F24").Value = 273
$ws.Range("F26").Value = 1094
$ws.Range("F29").Value = 917
$ws.Range("F30").Value = 892

$ws = $wb.Worksheets.Item("演出")
$ws.Range("F4").Value = 25
$ws.Range("F13").Value = 94
$ws.Range("F14").Value = 199
$ws.Range("F18").Value = 105
$ws.Range("F19").Value = 225
$ws.Range("F20").Value = 164
$ws.Range("F21").Value = 461

$ws = $wb.Worksheets.Item("本地生活")
$ws.Range("F6").Value = 467

$ws = $wb.Worksheets.Item("全部类型")
$ws.Range("F4").Value = 541
$ws.Range("F10").Value = 25
$ws.Range("F12").Value = 467
$ws.Range("F13").Value = 2397
$ws.Range("F14").Value = 71
$ws.Range("F24").Value = 1506
$ws.Range("F25").Value = 1506
$ws.Range("F28").Value = 740
$ws.Range("F29").Value = 94
$ws.Range("F30").Value = 199
$ws.Range("F31").Value = 1091
$ws.Range("F32").Value = 485
$ws.Range("F34").Value = 3404
$ws.Range("F35").Value = 634
$ws.Range("F36").Value = 3250
$ws.Range("F37").Value = 721
$ws.Range("F39").Value = 601
$ws.Range("F40").Value = 273
$ws.Range("F41").Value = 1094
$ws.Range("F43").Value = 105
$ws.Range("F44").Value = 225
$ws.Range("F45").Value = 164
$ws.Range("F46").Value = 461
$ws.Range("F49").Value = 917
$ws.Range("F50").Value = 892
